$d = $word.ActiveDocument

# The "FORM KELENGKAPAN ADMINISTRASI" checklist table is table #4 in the
# document. Row 5 of that table is point "3" -
# "Justifikasi / arahan penyusunan dokumen lingkungan" (fields
# {hasil_penapisan_yes}/{hasil_penapisan_no}/{hasil_penapisan_ket}) -
# which the commit message says to delete ("delete point in verification").
$t = $d.Tables.Item(4)
$t.Rows.Item(5).Delete()

# Every point after the deleted one shifts up, so its displayed number
# must be decremented by one (4->3, 5->4, ..., 11->10). After the delete,
# those rows are now at table rows 5..12, and each one's first cell still
# shows its old number.
for ($r = 5; $r -le 12; $r++) {
    $cell = $t.Rows.Item($r).Cells.Item(1)
    $oldNum = $r - 1
    $newNum = $r - 2
    $cell.Range.Text = [string]$newNum
}
